# Rename the worksheet "Property1" to "DataNode" as part of unifying the
# conception of DataNode / DataTable / Entity across the data config sheets.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"
